$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 is Idaho. The July 11, 2020 run succeeded for this row (it had
# previously timed out), so fill in the result columns with the new data
# and flip the status / "includes Hispanic Black" flag accordingly.

$ws.Range("B36").Value = "2020-07-10"
$ws.Range("B36").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("C36").Value = 9928
$ws.Range("D36").Value = 101
$ws.Range("E36").Value = 145
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 1.46
$ws.Range("H36").Value = 0.99

$ws.Range("J36").Value = $true

$ws.Range("O36").Value = "Success!"
